# Minor wording changes to #9 in MA2
#
# This script applies the content-visible text/formatting changes from the
# commit. (A number of hunks in the original diff are pure run-splitting /
# run-merging refactors with no visible text change -- e.g. " " + "("
# merging into " (", "R3" + " " merging into "R3 ", the hyperlink text
# runs recombining, and several paragraphs whose runs merge but whose
# concatenated text is identical. Those require no action here since the
# rendered text/formatting is unaffected.)

$d = $word.ActiveDocument
$apos = [char]0x2019

# ---------------------------------------------------------------------
# 1) "...contain two integer values in two's complement representation..."
#    -> insert " 8-bit" so it reads "...in 8-bit two's complement..."
# ---------------------------------------------------------------------
$old1 = "contain two integer values in two" + $apos + "s complement representation"
$new1 = "contain two integer values in 8-bit two" + $apos + "s complement representation"
$found1 = $d.Content.Find.Execute($old1, $true, $false, $false, $false, $false, $true, 1, $false, $new1, 2)
Write-Host "Step1 (8-bit #1): $found1"

# ---------------------------------------------------------------------
# 2) "...operation will also compute the logical OR..."
#    -> "...operation will always compute the logical OR..." with
#       "always" underlined
# ---------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute("operation will also compute the logical OR", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
Write-Host "Step2 find: $found2"
if ($found2) {
    $prefixLen = ("operation will ").Length
    $alsoStart = $rng2.Start + $prefixLen
    $alsoEnd = $alsoStart + ("also").Length
    $alsoRng = $d.Range($alsoStart, $alsoEnd)
    Write-Host "Step2 target text: [$($alsoRng.Text)]"
    $alsoRng.Text = "always"
    $alsoRng.Font.Underline = 1
}

# ---------------------------------------------------------------------
# 3) "...given two two's complement integer values in R0 and R1..."
#    -> insert "8-bit " so it reads "...given two 8-bit two's complement..."
# ---------------------------------------------------------------------
$old3 = "given two two" + $apos + "s complement integer values in R0 and R1"
$new3 = "given two 8-bit two" + $apos + "s complement integer values in R0 and R1"
$found3 = $d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)
Write-Host "Step3 (8-bit #2): $found3"

# ---------------------------------------------------------------------
# 4) Reword the "Give a counter example ..." sentence:
#    "Give a counter example (i.e. give binary values for R0 and R1)
#     that shows that the bitwise AND operation does not always compute
#     the logical AND operation."
#    ->
#    "Give a counter example showing that the bitwise AND operation does
#     not always compute the logical AND operation (i.e. give binary
#     values for R0 and R1 and the result of R0 & R1)."
#    with "does not always" underlined and R0/R1 (Courier) runs restored.
# ---------------------------------------------------------------------
$old4 = "Give a counter example (i.e. give binary values for R0 and R1) that shows that the bitwise AND operation does not always compute the logical AND operation."
$new4 = "Give a counter example showing that the bitwise AND operation does not always compute the logical AND operation (i.e. give binary values for R0 and R1 and the result of R0 & R1)."
$rng4 = $d.Content
$found4 = $rng4.Find.Execute($old4, $true, $false, $false, $false, $false, $true, 1, $false, $new4, 2)
Write-Host "Step4 replace: $found4"

if ($found4) {
    $rng4b = $d.Content
    $found4b = $rng4b.Find.Execute("Give a counter example showing", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    Write-Host "Step4 relocate: $found4b"
    if ($found4b) {
        $base = $rng4b.Start

        # "does not always" -> underline
        $u = $d.Range($base + 62, $base + 77)
        Write-Host "Step4 underline text: [$($u.Text)]"
        $u.Font.Underline = 1

        # first "R0" (in "give binary values for R0 and R1") -> Courier
        $c1 = $d.Range($base + 141, $base + 143)
        Write-Host "Step4 courier1 text: [$($c1.Text)]"
        $c1.Font.Name = "Courier"

        # first "R1" -> Courier
        $c2 = $d.Range($base + 148, $base + 150)
        Write-Host "Step4 courier2 text: [$($c2.Text)]"
        $c2.Font.Name = "Courier"

        # "R0 &" -> Courier
        $c3 = $d.Range($base + 169, $base + 173)
        Write-Host "Step4 courier3 text: [$($c3.Text)]"
        $c3.Font.Name = "Courier"

        # " R1" -> Courier
        $c4 = $d.Range($base + 173, $base + 176)
        Write-Host "Step4 courier4 text: [$($c4.Text)]"
        $c4.Font.Name = "Courier"
    }
}

Write-Host "Done"
